# Conserto do erro com o rótulo da coluna 2050 nas tabelas e
# retirada das linhas com total das tabelas.

$wb = $excel.ActiveWorkbook

# Sheets 1-4 share the same layout: header row 1, data rows 2-12,
# a "Total" row 13, and column E labeled with the last period.
# Sheet 4 uses a year-range label ("2031-2040" -> "2041-2050"),
# the others use a single year ("2040" -> "2050").
$rangeSheetLabels = @{ 1 = "2050"; 2 = "2050"; 3 = "2050"; 4 = "2041-2050" }
foreach ($idx in $rangeSheetLabels.Keys) {
    $ws = $wb.Worksheets.Item($idx)

    $ws.Range("E1").Value = $rangeSheetLabels[$idx]

    # Remove the "Total" row (row 13) entirely.
    $ws.Rows.Item(13).Delete()
}

# Sheet 5 only needs the column header fix; it has no Total row.
$ws5 = $wb.Worksheets.Item(5)
$ws5.Range("E1").Value = "2050"

# Sheet 6 only needs the "Total" row (row 4) removed; it has no
# column E / 2050 label at all.
$ws6 = $wb.Worksheets.Item(6)
$ws6.Rows.Item(4).Delete()
